# Applies the cryptos-list refresh described by the commit:
# "Updated cryptos list on Mon Jul 31 03:08:31 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.455.51"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.41"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7076"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.73"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3161"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07876"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.72"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08000"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -3.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.889.32"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.217"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.13"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7057"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.518"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.495.69"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008368"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -3.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "257.15"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +5.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.133.27"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.23"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.638"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1560"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.084"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.11"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.88"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.502"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.342"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -2.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.262"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.210"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.899"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7505"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.175"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01882"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.266.41"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.756"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9022"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.60"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.983"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -8.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.83"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.030.74"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.556"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4337"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.16%  "

Write-Host "Applied cryptos update"
